# Applies the cryptos.xlsx price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.025.91"
$ws.Range("E2").Value = "  -3.24%  "
$ws.Range("D3").Value = "3.236.82"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.49%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.235.84"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.459"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("E11").Value = "  -5.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.59%  "
$ws.Range("D13").Value = "3.785.46"
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.54%  "
$ws.Range("D16").Value = "3.229.21"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").Value = "59.123.78"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.19%  "
$ws.Range("E20").Value = "  -6.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "362.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.520"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.63%  "
$ws.Range("D26").Value = "3.360.44"
$ws.Range("E26").Value = "  -4.69%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0₃0978"
$ws.Range("E27").Value = "  -10.36%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.171"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.20%  "
$ws.Range("E33").Value = "  -8.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("E35").Value = "  -3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.63%  "
$ws.Range("E39").Value = "  -6.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0709"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("D42").Value = "3.265.03"
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.718"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.65%  "
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.59%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "2.296.29"
$ws.Range("E49").Value = "  -8.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.04%  "
